# datannur "info.xlsx" live-update stamp
#
# The source sheet holds two rows of metadata in column A/B:
#   A1: id        B1: value
#   A2: last_update   B2: <unix timestamp>
#
# This edit bumps the stored timestamp (footer "last update" info that is
# live-updated on every export), auto-sizes column B to fit the new value,
# and leaves the active selection on D8 the way the authoring app did when
# it saved the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "last_update" timestamp stored in B2.
$ws.Range("B2").Value = 1728846186

# Column B is re-sized ("best fit") to the new value so the number is not
# truncated, matching the <cols> width tweak introduced by this commit
# (Excel auto-fit of column B to the 10-digit timestamp -> ~11.16 chars).
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(2).ColumnWidth = 10.3

# The author's last selection when the file was saved was D8.
$ws.Range("D8").Select() | Out-Null
